# Edit script: adds "PO Forecast" sheet, renames header cells on existing sheets,
# and populates forecast data, matching the target diff.

$wb = $excel.ActiveWorkbook

# --- 1. Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "PO Forecast"

# Match page margins used by the other sheets (0.75in/0.75in/1in/1in/0.5in/0.5in)
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# --- 3. Apply styles by copying from existing styled cells (reuses the same style ids) ---
# Header style (bold, centered, bordered) comes from the "Requested quantity" header cell.
$wsWeekly.Range("B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Date style (custom date/time number format) comes from column A data cells.
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A70").PasteSpecial(-4122)

# --- 4. Write header text ---
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# --- 5. Write forecast data rows (row, ds, PO_Forecast, yhat_lower, yhat_upper) ---
$forecastData = @"
2	44948.99999999999	413	70.43168026788176	779.8622815165043
3	44976.99999999999	402	18.09212415780437	755.5063702797789
4	44983.99999999999	399	43.69328829306145	750.8460508065443
5	44990.99999999999	396	24.18720275955274	739.5646792250959
6	44997.99999999999	393	42.81293576507756	762.1153529541472
7	45011.99999999999	387	16.79322990289989	759.450020376682
8	45025.99999999999	382	29.93709935290198	745.7175967843556
9	45088.99999999999	356	11.90789125252416	729.3343620760595
10	45095.99999999999	353	-18.34716961785314	718.8211067961756
11	45102.99999999999	351	-34.93104116253176	730.0984114510791
12	45109.99999999999	348	-19.93936015467605	691.2304999076084
13	45116.99999999999	345	8.21486858464888	745.0358208765069
14	45123.99999999999	342	-26.23242665485563	710.301785073809
15	45130.99999999999	339	-17.21870631839098	677.4974071411968
16	45137.99999999999	336	-9.434664188397294	679.3922240567549
17	45151.99999999999	331	-39.12562761336237	692.8322487460122
18	45165.99999999999	325	-65.44269417990544	703.568699142292
19	45179.99999999999	320	-43.24331724522975	693.1749904497657
20	45186.99999999999	317	-32.26587177523196	677.8939816842723
21	45193.99999999999	314	-32.82946223105802	651.342631821902
22	45200.99999999999	311	-66.62575202645844	691.0273669271616
23	45221.99999999999	303	-42.64520269474779	657.6394549592487
24	45228.99999999999	300	-37.20819199084197	675.6893371927912
25	45235.99999999999	297	-53.35920286094024	654.6620059540086
26	45242.99999999999	294	-45.9642407616453	636.9959179593583
27	45249.99999999999	291	-104.0996904082291	671.2442641492967
28	45256.99999999999	288	-54.35097602843091	668.2863113583636
29	45270.99999999999	283	-82.7509675045762	630.3337288756024
30	45277.99999999999	280	-67.49881272099242	627.4129943831256
31	45298.99999999999	271	-62.86395098234825	631.4906131177272
32	45312.99999999999	266	-72.47048015592652	649.0445033346593
33	45319.99999999999	263	-119.9785964680939	623.4471750998864
34	45326.99999999999	260	-88.08262175102611	645.9126451791257
35	45333.99999999999	257	-80.63436416651722	612.9143494956437
36	45340.99999999999	254	-103.1481939206104	631.0517723483969
37	45347.99999999999	252	-102.7264591283892	592.5550917738354
38	45354.99999999999	249	-108.1131830075013	578.2048644137999
39	45361.99999999999	246	-109.5309432832492	583.277036469898
40	45368.99999999999	243	-118.7301595721049	622.3711163658712
41	45375.99999999999	240	-114.4900595838692	582.0914604544917
42	45382.99999999999	238	-110.696363319089	594.1817841017527
43	45389.99999999999	235	-153.689085137783	593.2547798609691
44	45396.99999999999	232	-132.9260071330723	590.2719711837428
45	45403.99999999999	229	-119.7312544440371	583.0540671519844
46	45424.99999999999	221	-126.6131553346258	573.56387092229
47	45431.99999999999	218	-157.9937898631366	561.2214987408386
48	45445.99999999999	212	-121.0649514618336	585.1830932146308
49	45459.99999999999	206	-150.8455065643812	576.3733565610969
50	45466.99999999999	204	-149.400271081796	543.1976194155135
51	45473.99999999999	201	-156.1833096879146	585.1931201528608
52	45480.99999999999	198	-142.6289844810346	538.8077823708609
53	45487.99999999999	195	-147.2006552786007	562.677002322769
54	45501.99999999999	189	-171.5162386035038	537.9179119578868
55	45515.99999999999	184	-188.4194515076801	524.4009892747946
56	45557.99999999999	167	-187.5242186221897	521.0679957213746
57	45564.99999999999	164	-180.723198664078	535.8828645681903
58	45571.99999999999	161	-233.0961780375264	522.1865242789758
59	45578.99999999999	158	-211.0485354308513	518.4484754355709
60	45585.99999999999	155	-198.9305991450502	505.8845531811442
61	45592.99999999999	153	-198.7979964492282	489.8020361024155
62	45599.99999999999	150	-185.2519354151383	511.5453413017725
63	45606.99999999999	147	-217.4210045132811	509.3488508930305
64	45613.99999999999	144	-204.2053242120054	519.8800334748614
65	45620.99999999999	141	-210.2320848470873	506.6930490888126
66	45627.99999999999	139	-224.4268435573545	482.2088871124491
67	45634.99999999999	136	-240.2628759877772	503.7657649484183
68	45641.99999999999	133	-248.4792157509626	492.269744223363
69	45648.99999999999	130	-225.0374902790714	495.3807748892797
70	45655.99999999999	127	-244.6286949623754	475.3064783621427
"@

$rows = $forecastData -split "`n"
foreach ($line in $rows) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split "`t"
    $r = [int]$parts[0]
    $dsVal = [double]$parts[1]
    $poVal = [double]$parts[2]
    $lowerVal = [double]$parts[3]
    $upperVal = [double]$parts[4]

    $newSheet.Cells.Item($r, 1).Value = $dsVal
    $newSheet.Cells.Item($r, 2).Value = $poVal
    $newSheet.Cells.Item($r, 3).Value = $lowerVal
    $newSheet.Cells.Item($r, 4).Value = $upperVal
}
